# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "K" values for column G (rows 2-15), replacing the old Strike# based values.
$kValues = @{
    2  = 2
    3  = 2
    4  = 1
    5  = 2
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 2
    11 = 3
    12 = 1
    13 = 0
    14 = 1
    15 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
